$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.005899579586891775
$ws.Range("J2").Value = 0.008823342375055644
$ws.Range("M2").Value = 2.6796385
$ws.Range("N2").Value = 5.359277000000001
$ws.Range("O2").Value = 0.03934793987875059
$ws.Range("P2").Value = 0.02842274181890429
$ws.Range("Q2").Value = 0.1722650270366667
$ws.Range("R2").Value = 1.03359016222
$ws.Range("S2").Value = 0.0002321363028949218
$ws.Range("T2").Value = 0.0002507835823060044
$ws.Range("I3").Value = 0.005899579586891775
$ws.Range("J3").Value = 0.008823342375055644
$ws.Range("M3").Value = 50.102415
$ws.Range("O3").Value = 0.7357062578404556
$ws.Range("P3").Value = 0.7971493203553003
$ws.Range("Q3").Value = 3.2209172523
$ws.Range("R3").Value = 28.9882552707
$ws.Range("S3").Value = 0.004340357620704089
$ws.Range("T3").Value = 0.007033521377537728
$ws.Range("I4").Value = 0.005899579586891775
$ws.Range("J4").Value = 0.008823342375055644
$ws.Range("M4").Value = 0.2784063333333333
$ws.Range("N4").Value = 0.8352189999999999
$ws.Range("O4").Value = 0.004088131912518571
$ws.Range("P4").Value = 0.00442955532980352
$ws.Range("Q4").Value = 0.01789781514888889
$ws.Range("R4").Value = 0.16108033634
$ws.Range("S4").Value = 0.00002411825957961539
$ws.Range("T4").Value = 0.00003908348324410898
$ws.Range("I5").Value = 0.005899579586891775
$ws.Range("J5").Value = 0.008823342375055644
$ws.Range("M5").Value = 13.0677535
$ws.Range("N5").Value = 26.135507
$ws.Range("O5").Value = 0.1918875173156127
$ws.Range("P5").Value = 0.1386087652806835
$ws.Range("Q5").Value = 0.8400823133366666
$ws.Range("R5").Value = 5.04049388002
$ws.Range("S5").Value = 0.001132055680134531
$ws.Range("T5").Value = 0.001222992592255196
$ws.Range("I6").Value = 0.005899579586891775
$ws.Range("J6").Value = 0.008823342375055644
$ws.Range("M6").Value = 1.863198333333333
$ws.Range("N6").Value = 5.589594999999999
$ws.Range("O6").Value = 0.02735929342789644
$ws.Range("P6").Value = 0.02964422543511714
$ws.Range("Q6").Value = 0.1197788101888889
$ws.Range("R6").Value = 1.0780092917
$ws.Range("S6").Value = 0.0001614083290190002
$ws.Range("T6").Value = 0.0002615611504573714
$ws.Range("I7").Value = 0.005899579586891775
$ws.Range("J7").Value = 0.008823342375055644
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1097013333333333
$ws.Range("N7").Value = 0.329104
$ws.Range("O7").Value = 0.001610859624766094
$ws.Range("P7").Value = 0.001745391780191372
$ws.Range("Q7").Value = 0.00705233304888889
$ws.Range("R7").Value = 0.06347099744
$ws.Range("S7").Value = 0.00000950339455961819
$ws.Range("T7").Value = 0.00001540018925523634
$ws.Range("G8").Value = 10.832535
$ws.Range("H8").Value = 21.66507
$ws.Range("I8").Value = 0.9941004204131083
$ws.Range("J8").Value = 0.9911766576249443
$ws.Range("M8").Value = 2.6796385
$ws.Range("N8").Value = 5.359277000000001
$ws.Range("O8").Value = 0.03934793987875059
$ws.Range("P8").Value = 0.02842274181890429
$ws.Range("Q8").Value = 29.0272778385975
$ws.Range("R8").Value = 116.10911135439
$ws.Range("S8").Value = 0.03911580357585567
$ws.Range("T8").Value = 0.02817195823659829
$ws.Range("G9").Value = 10.832535
$ws.Range("H9").Value = 21.66507
$ws.Range("I9").Value = 0.9941004204131083
$ws.Range("J9").Value = 0.9911766576249443
$ws.Range("M9").Value = 50.102415
$ws.Range("O9").Value = 0.7357062578404556
$ws.Range("P9").Value = 0.7971493203553003
$ws.Range("Q9").Value = 542.736164072025
$ws.Range("R9").Value = 3256.41698443215
$ws.Range("S9").Value = 0.7313659002197516
$ws.Range("T9").Value = 0.7901157989777625
$ws.Range("G10").Value = 10.832535
$ws.Range("H10").Value = 21.66507
$ws.Range("I10").Value = 0.9941004204131083
$ws.Range("J10").Value = 0.9911766576249443
$ws.Range("M10").Value = 0.2784063333333333
$ws.Range("N10").Value = 0.8352189999999999
$ws.Range("O10").Value = 0.004088131912518571
$ws.Range("P10").Value = 0.00442955532980352
$ws.Range("Q10").Value = 3.015846350055
$ws.Range("R10").Value = 18.09507810033
$ws.Range("S10").Value = 0.004064013652938956
$ws.Range("T10").Value = 0.004390471846559411
$ws.Range("G11").Value = 10.832535
$ws.Range("H11").Value = 21.66507
$ws.Range("I11").Value = 0.9941004204131083
$ws.Range("J11").Value = 0.9911766576249443
$ws.Range("M11").Value = 13.0677535
$ws.Range("N11").Value = 26.135507
$ws.Range("O11").Value = 0.1918875173156127
$ws.Range("P11").Value = 0.1386087652806835
$ws.Range("Q11").Value = 141.5568971601225
$ws.Range("R11").Value = 566.2275886404899
$ws.Range("S11").Value = 0.1907554616354782
$ws.Range("T11").Value = 0.1373857726884283
$ws.Range("G12").Value = 10.832535
$ws.Range("H12").Value = 21.66507
$ws.Range("I12").Value = 0.9941004204131083
$ws.Range("J12").Value = 0.9911766576249443
$ws.Range("M12").Value = 1.863198333333333
$ws.Range("N12").Value = 5.589594999999999
$ws.Range("O12").Value = 0.02735929342789644
$ws.Range("P12").Value = 0.02964422543511714
$ws.Range("Q12").Value = 20.183161157775
$ws.Range("R12").Value = 121.09896694665
$ws.Range("S12").Value = 0.02719788509887745
$ws.Range("T12").Value = 0.02938266428465977
$ws.Range("G13").Value = 10.832535
$ws.Range("H13").Value = 21.66507
$ws.Range("I13").Value = 0.9941004204131083
$ws.Range("J13").Value = 0.9911766576249443
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1097013333333333
$ws.Range("N13").Value = 0.329104
$ws.Range("O13").Value = 0.001610859624766094
$ws.Range("P13").Value = 0.001745391780191372
$ws.Range("Q13").Value = 1.18834353288
$ws.Range("R13").Value = 7.13006119728
$ws.Range("S13").Value = 0.001601356230206476
$ws.Range("T13").Value = 0.001729991590936136
